$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.354.34'
$ws.Range('E2').Value = '  -1.03%  '
$ws.Range('D3').Value = '3.246.05'
$ws.Range('E3').Value = '  +3.07%  '
$ws.Range('E4').Value = '  +0.03%  '
$c = $ws.Range('D5')
$origStyle = $c.Style()
$c.NumberFormat = '@'
$c.Value = '594.94'
$c.Style = $origStyle
$ws.Range('E5').Value = '  -1.17%  '
$c = $ws.Range('D6')
$origStyle = $c.Style()
$c.NumberFormat = '@'
$c.Value = '140.47'
$c.Style = $origStyle
$ws.Range('E6').Value = '  -1.13%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '3.239.04'
$ws.Range('E8').Value = '  +3.01%  '
$ws.Range('E9').Value = '  -1.75%  '
$ws.Range('E10').Value = '  -1.03%  '
$ws.Range('E11').Value = '  -0.63%  '
$c = $ws.Range('D12')
$origStyle = $c.Style()
$c.NumberFormat = '@'
$c.Value = '0.464'
$c.Style = $origStyle
$ws.Range('E12').Value = '  -0.56%  '
$ws.Range('E13').Value = '  -3.00%  '
$c = $ws.Range('D14')
$origStyle = $c.Style()
$c.NumberFormat = '@'
$c.Value = '34.30'
$c.Style = $origStyle
$ws.Range('E14').Value = '  -1.70%  '
$ws.Range('D15').Value = '3.771.81'
$ws.Range('E15').Value = '  +2.92%  '
$ws.Range('E16').Value = '  -0.38%  '
$ws.Range('D17').Value = '3.245.07'
$ws.Range('E17').Value = '  +3.24%  '
$ws.Range('D18').Value = '63.341.71'
$ws.Range('E18').Value = '  -1.01%  '
$ws.Range('E19').Value = '  -1.14%  '
$c = $ws.Range('D20')
$origStyle = $c.Style()
$c.NumberFormat = '@'
$c.Value = '473.72'
$c.Style = $origStyle
$c = $ws.Range('D21')
$origStyle = $c.Style()
$c.NumberFormat = '@'
$c.Value = '14.18'
$c.Style = $origStyle
$ws.Range('E21').Value = '  -3.39%  '
$c = $ws.Range('D22')
$origStyle = $c.Style()
$c.NumberFormat = '@'
$c.Value = '0.731'
$c.Style = $origStyle
$ws.Range('E22').Value = '  +2.74%  '
$ws.Range('E23').Value = '  +2.59%  '
$c = $ws.Range('D24')
$origStyle = $c.Style()
$c.NumberFormat = '@'
$c.Value = '84.14'
$c.Style = $origStyle
$ws.Range('E24').Value = '  -4.88%  '
$c = $ws.Range('D25')
$origStyle = $c.Style()
$c.NumberFormat = '@'
$c.Value = '13.16'
$c.Style = $origStyle
$ws.Range('E25').Value = '  -0.48%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('E27').Value = '  -0.98%  '
$c = $ws.Range('D28')
$origStyle = $c.Style()
$c.NumberFormat = '@'
$c.Value = '7.30'
$c.Style = $origStyle
$ws.Range('E28').Value = '  +4.65%  '
$ws.Range('E29').Value = '  -1.31%  '
$ws.Range('E30').Value = '  +2.72%  '
$c = $ws.Range('D31')
$origStyle = $c.Style()
$c.NumberFormat = '@'
$c.Value = '27.52'
$c.Style = $origStyle
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('E32').Value = '  +0.00%  '
$ws.Range('E33').Value = '  -4.32%  '
$ws.Range('E34').Value = '  -4.60%  '
$c = $ws.Range('D35')
$origStyle = $c.Style()
$c.NumberFormat = '@'
$c.Value = '1.08'
$c.Style = $origStyle
$ws.Range('E35').Value = '  -1.65%  '
$ws.Range('E36').Value = '  -2.19%  '
$c = $ws.Range('D37')
$origStyle = $c.Style()
$c.NumberFormat = '@'
$c.Value = '52.68'
$c.Style = $origStyle
$ws.Range('E37').Value = '  -0.17%  '
$ws.Range('E38').Value = '  -4.37%  '
$ws.Range('E39').Value = '  -1.29%  '
$c = $ws.Range('D40')
$origStyle = $c.Style()
$c.NumberFormat = '@'
$c.Value = '422.65'
$c.Style = $origStyle
$ws.Range('E40').Value = '  -2.22%  '
$ws.Range('E41').Value = '  +0.11%  '
$ws.Range('D42').Value = '2.979.10'
$ws.Range('E42').Value = '  +2.36%  '
$c = $ws.Range('D43')
$origStyle = $c.Style()
$c.NumberFormat = '@'
$c.Value = '2.74'
$c.Style = $origStyle
$ws.Range('E43').Value = '  -6.22%  '
$ws.Range('E44').Value = '  -7.93%  '
$ws.Range('E45').Value = '  +2.15%  '
$ws.Range('E46').Value = '  -0.82%  '
$ws.Range('E47').Value = '  +0.07%  '
$c = $ws.Range('D48')
$origStyle = $c.Style()
$c.NumberFormat = '@'
$c.Value = '25.89'
$c.Style = $origStyle
$ws.Range('E48').Value = '  +0.44%  '
$ws.Range('E49').Value = '  -3.48%  '
$ws.Range('E50').Value = '  -0.57%  '
$c = $ws.Range('D51')
$origStyle = $c.Style()
$c.NumberFormat = '@'
$c.Value = '121.41'
$c.Style = $origStyle
$ws.Range('E51').Value = '  +0.43%  '
